# Masterfile_Heraclee.xlsx - add "NEWS PAGES" + "AUTHORISATION" translation
# blocks (KEYS / FR / EN) to Sheet1, matching the "translations modified
# masterfile for news" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (password 9488) - unprotect so the new rows can
# be written; the source workbook no longer protects this sheet either.
$ws.Unprotect("9488")

# --- 1. write the new cell values --------------------------------------
# Values are written in the same order the original author typed them in
# (KEYS column first, then the EN column, then the FR column) so that the
# shared-string table gets rebuilt with the same ordering as the target.

# KEYS (column A) -- //NEWS PAGES block first, then latest_news row
$ws.Cells.Item(162,1).Value = '//NEWS PAGES'
$ws.Cells.Item(163,1).Value = 'latest_news'
$ws.Cells.Item(163,3).Value = 'Latest news'
$ws.Cells.Item(163,2).Value = 'Dernières nouvelles'

# // AUTHORISATION block + rest of the KEYS column
$ws.Cells.Item(157,1).Value = '// AUTHORISATION'
$ws.Cells.Item(158,1).Value = 'password'
$ws.Cells.Item(159,1).Value = 'remember_me'
$ws.Cells.Item(160,1).Value = 'log_in'
$ws.Cells.Item(164,1).Value = 'list_of_articles'
$ws.Cells.Item(165,1).Value = 'published'
$ws.Cells.Item(166,1).Value = 'title'
$ws.Cells.Item(167,1).Value = 'posting_date'
$ws.Cells.Item(168,1).Value = 'article_text'
$ws.Cells.Item(169,1).Value = 'save'
$ws.Cells.Item(170,1).Value = 'publish'
$ws.Cells.Item(171,1).Value = 'creating_article'
$ws.Cells.Item(172,1).Value = 'editing_article'
$ws.Cells.Item(173,1).Value = 'choose_header_img'
$ws.Cells.Item(174,1).Value = 'choose_body_img'

# EN (column C)
$ws.Cells.Item(158,3).Value = 'Password'
$ws.Cells.Item(159,3).Value = 'Remember me'
$ws.Cells.Item(160,3).Value = 'Log in'
$ws.Cells.Item(164,3).Value = 'List of articles'
$ws.Cells.Item(165,3).Value = 'Published'
$ws.Cells.Item(166,3).Value = 'Title'
$ws.Cells.Item(167,3).Value = 'Posting date'
$ws.Cells.Item(168,3).Value = 'Article text'
$ws.Cells.Item(169,3).Value = 'Save'
$ws.Cells.Item(170,3).Value = 'Publish'
$ws.Cells.Item(171,3).Value = 'Creating new article'
$ws.Cells.Item(172,3).Value = 'Editing article'
$ws.Cells.Item(173,3).Value = 'Choose a header image'
$ws.Cells.Item(174,3).Value = 'Choose a body image'

# FR (column B)
$ws.Cells.Item(159,2).Value = 'Se souvenir de moi'
$ws.Cells.Item(158,2).Value = 'Mot de passe'
$ws.Cells.Item(160,2).Value = 'Connexion'
$ws.Cells.Item(164,2).Value = 'Liste des articles'
$ws.Cells.Item(165,2).Value = 'Publié'
$ws.Cells.Item(166,2).Value = 'Titre'
$ws.Cells.Item(167,2).Value = 'Date de publication'
$ws.Cells.Item(168,2).Value = 'Texte de l''article'
$ws.Cells.Item(169,2).Value = 'Enregistrer'
$ws.Cells.Item(170,2).Value = 'Publier'
$ws.Cells.Item(171,2).Value = 'Création d’article'
$ws.Cells.Item(172,2).Value = 'Modification de l''article'
$ws.Cells.Item(173,2).Value = 'Choisissez une image d''en-tête'
$ws.Cells.Item(174,2).Value = 'Choisissez une image du corps'

# row 156 is left blank on purpose (spacer row before "// AUTHORISATION"),
# as is row 161 (spacer row before "//NEWS PAGES") which gets no cells at
# all.

# --- 2. formatting --------------------------------------------------------
# Column A of the new block reuses the section-header look already used by
# A154/A155 ("I am interested" row); columns B/C reuse the wrapped body
# look already used by B154:C155.
$ws.Range("A154").Copy()
$ws.Range("A156:A160").PasteSpecial(-4122)
$ws.Range("A154").Copy()
$ws.Range("A162:A174").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B154:C154").Copy()
$ws.Range("B158:C160").PasteSpecial(-4122)
$ws.Range("B154:C154").Copy()
$ws.Range("B163:C174").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. view / selection state ------------------------------------------
$ws.Range("B175").Select()
$excel.ActiveWindow.ScrollRow = 145
$excel.ActiveWindow.ScrollColumn = 1

$wb.Windows.Item(1).Left = 34620
$wb.Windows.Item(1).Top = 520

# --- 4. leave the sheet unprotected, as in the target workbook -----------
# (Unprotect above already removes <sheetProtection> from the saved file.)
